{"js": "// Replace the 25 three-digit \u00f7 one-digit division problem/answer strings\n// in the table cells with their new values, preserving all formatting.\nconst replacements = [\n  [\"180\u00f77=25, 5\", \"463\u00f75=92, 3\"],\n  [\"349\u00f79=38, 7\", \"829\u00f78=103, 5\"],\n  [\"714\u00f73=238, 0\", \"579\u00f76=96, 3\"],\n  [\"593\u00f79=65, 8\", \"702\u00f74=175, 2\"],\n  [\"700\u00f73=233, 1\", \"822\u00f77=117, 3\"],\n  [\"332\u00f76=55, 2\", \"503\u00f75=100, 3\"],\n  [\"634\u00f76=105, 4\", \"854\u00f79=94, 8\"],\n  [\"645\u00f78=80, 5\", \"231\u00f73=77, 0\"],\n  [\"771\u00f73=257, 0\", \"375\u00f79=41, 6\"],\n  [\"936\u00f79=104, 0\", \"712\u00f73=237, 1\"],\n  [\"512\u00f77=73, 1\", \"657\u00f75=131, 2\"],\n  [\"239\u00f72=119, 1\", \"840\u00f74=210, 0\"],\n  [\"517\u00f72=258, 1\", \"658\u00f78=82, 2\"],\n  [\"610\u00f74=152, 2\", \"323\u00f72=161, 1\"],\n  [\"635\u00f79=70, 5\", \"915\u00f77=130, 5\"],\n  [\"108\u00f74=27, 0\", \"455\u00f75=91, 0\"],\n  [\"822\u00f73=274, 0\", \"232\u00f77=33, 1\"],\n  [\"987\u00f79=109, 6\", \"470\u00f77=67, 1\"],\n  [\"612\u00f79=68, 0\", \"463\u00f79=51, 4\"],\n  [\"279\u00f79=31, 0\", \"115\u00f77=16, 3\"],\n  [\"521\u00f76=86, 5\", \"305\u00f79=33, 8\"],\n  [\"456\u00f79=50, 6\", \"724\u00f77=103, 3\"],\n  [\"180\u00f75=36, 0\", \"351\u00f77=50, 1\"],\n  [\"324\u00f75=64, 4\", \"732\u00f76=122, 0\"],\n  [\"602\u00f76=100, 2\", \"328\u00f74=82, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 three-digit divided-by-one-digit division problem/answer\n# strings in the table cells with their new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"180\u00f77=25, 5\", \"463\u00f75=92, 3\"),\n  @(\"349\u00f79=38, 7\", \"829\u00f78=103, 5\"),\n  @(\"714\u00f73=238, 0\", \"579\u00f76=96, 3\"),\n  @(\"593\u00f79=65, 8\", \"702\u00f74=175, 2\"),\n  @(\"700\u00f73=233, 1\", \"822\u00f77=117, 3\"),\n  @(\"332\u00f76=55, 2\", \"503\u00f75=100, 3\"),\n  @(\"634\u00f76=105, 4\", \"854\u00f79=94, 8\"),\n  @(\"645\u00f78=80, 5\", \"231\u00f73=77, 0\"),\n  @(\"771\u00f73=257, 0\", \"375\u00f79=41, 6\"),\n  @(\"936\u00f79=104, 0\", \"712\u00f73=237, 1\"),\n  @(\"512\u00f77=73, 1\", \"657\u00f75=131, 2\"),\n  @(\"239\u00f72=119, 1\", \"840\u00f74=210, 0\"),\n  @(\"517\u00f72=258, 1\", \"658\u00f78=82, 2\"),\n  @(\"610\u00f74=152, 2\", \"323\u00f72=161, 1\"),\n  @(\"635\u00f79=70, 5\", \"915\u00f77=130, 5\"),\n  @(\"108\u00f74=27, 0\", \"455\u00f75=91, 0\"),\n  @(\"822\u00f73=274, 0\", \"232\u00f77=33, 1\"),\n  @(\"987\u00f79=109, 6\", \"470\u00f77=67, 1\"),\n  @(\"612\u00f79=68, 0\", \"463\u00f79=51, 4\"),\n  @(\"279\u00f79=31, 0\", \"115\u00f77=16, 3\"),\n  @(\"521\u00f76=86, 5\", \"305\u00f79=33, 8\"),\n  @(\"456\u00f79=50, 6\", \"724\u00f77=103, 3\"),\n  @(\"180\u00f75=36, 0\", \"351\u00f77=50, 1\"),\n  @(\"324\u00f75=64, 4\", \"732\u00f76=122, 0\"),\n  @(\"602\u00f76=100, 2\", \"328\u00f74=82, 0\"),\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Execute(\n    $oldText,   # FindText\n    $false,     # MatchCase\n    $true,      # MatchWholeWord\n    $false,     # MatchWildcards\n    $false,     # MatchSoundsLike\n    $false,     # MatchAllWordForms\n    $true,      # Forward\n    1,          # Wrap (wdFindContinue)\n    $false,     # Format\n    $newText,   # ReplaceWith\n    2           # Replace (wdReplaceAll)\n  ) | Out-Null\n}\n"}
